# Diary table update
# Week "14th Nov" row = table row 4
# Week "21st Nov" row = table row 5
# Week "28th Nov" row = table row 6
# Columns: 1=Date 2=Mon 3=Tue 4=Wed 5=Thu 6=Fri 7=Sat 8=Sun

$d = $word.ActiveDocument
$t = $d.Tables.Item(1)
$wns = "xmlns:w='http://schemas.openxmlformats.org/wordprocessingml/2006/main'"

function Set-CellXml($cell, $innerParasXml) {
    $cell.Range.InsertXML($innerParasXml)
}

# ---- Row 4 (14th Nov) ----

# Tue cell: keep the same text but drop the _GoBack bookmark (it moves to the
# end of the document, in the new final row-6 cell, below).
$c = $t.Cell(4, 3)
$xml = "<w:p $wns><w:r><w:t>Produced timescale document, to be discussed with team Weds.</w:t></w:r></w:p>"
Set-CellXml $c $xml

# Wed cell: two new paragraphs
$c = $t.Cell(4, 4)
$xml = "<w:p $wns><w:r><w:t>Attended team meeting</w:t></w:r></w:p>" +
       "<w:p $wns><w:r><w:t>Introduced team to timescale document, made final changes to it.</w:t></w:r></w:p>"
Set-CellXml $c $xml

# Thu cell (col 5) is left untouched (stays empty)

# Fri cell: one new paragraph
$c = $t.Cell(4, 6)
$xml = "<w:p $wns><w:r><w:t>Worked with Tassos on database creation, basic project layout and  interfacing with database.</w:t></w:r></w:p>"
Set-CellXml $c $xml

# ---- Row 5 (21st Nov) ----

# Mon cell: two new paragraphs
$c = $t.Cell(5, 2)
$xml = "<w:p $wns><w:r><w:t>Attended team meeting</w:t></w:r></w:p>" +
       "<w:p $wns><w:r><w:t>Worked on UML</w:t></w:r></w:p>"
Set-CellXml $c $xml

# Tue cell (col 3) stays empty

# Wed cell: one new paragraph
$c = $t.Cell(5, 4)
$xml = "<w:p $wns><w:r><w:t>Attended team meeting</w:t></w:r></w:p>"
Set-CellXml $c $xml

# Thu, Fri cells (col 5, 6) stay empty

# Sat cell: three runs in one paragraph, last run preceded by a rendered page break marker
$c = $t.Cell(5, 7)
$xml = "<w:p $wns>" +
       "<w:r><w:t>Worked on Test Plan</w:t></w:r>" +
       "<w:r><w:t xml:space=`"preserve`"> and </w:t></w:r>" +
       "<w:r><w:lastRenderedPageBreak/><w:t>researched Java API frameworks</w:t></w:r>" +
       "</w:p>"
Set-CellXml $c $xml

# Sun cell: two runs in one paragraph, each preceded by a rendered page break marker
$c = $t.Cell(5, 8)
$xml = "<w:p $wns>" +
       "<w:r><w:lastRenderedPageBreak/><w:t xml:space=`"preserve`">Worked on Test Plan and </w:t></w:r>" +
       "<w:r><w:lastRenderedPageBreak/><w:t>researched Java API frameworks</w:t></w:r>" +
       "</w:p>"
Set-CellXml $c $xml

# ---- Row 6 (28th Nov) ----

# Mon cell: four new paragraphs, the _GoBack bookmark now sits at the end of
# the last paragraph (it was removed from row 4 / Tue above).
$c = $t.Cell(6, 2)
$xml = "<w:p $wns><w:r><w:t>Attended team meeting</w:t></w:r></w:p>" +
       "<w:p $wns><w:r><w:t>Arranged next Sprint</w:t></w:r></w:p>" +
       "<w:p $wns><w:r><w:t>Had meeting with Julian</w:t></w:r></w:p>" +
       "<w:p $wns><w:r><w:t>Tested Java API framework</w:t></w:r>" +
       "<w:bookmarkStart w:id=`"0`" w:name=`"_GoBack`"/><w:bookmarkEnd w:id=`"0`"/></w:p>"
Set-CellXml $c $xml

Write-Output "Diary table updated"
